$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet currently has (rows 1-16):
#   Row 5: SetReportErrCode / F5 blank
#   Row 6: is_fail (styled s2 on A) / height 35.65
#   Rows 7-16: blank filler rows (D:G only)
#
# Target layout (rows 1-18):
#   Row 5: SetReportErrCode, F5 gets new docstring, height 24
#   Row 6: NEW "ReportError" row, height 47.25
#   Row 7: NEW "ResetWarning" row, height 24
#   Row 8: is_fail (moved down from old row 6), height 35.65
#   Rows 9-16: blank filler rows (shifted down from old rows 7-14)
#   Rows 17-18: NEW blank filler rows
#
# Strategy: insert two blank rows at 8:9 (both neighbouring rows 7 & 8
# are already unstyled filler rows there) so Excel does not have to
# fabricate brand-new style records for the insert - this also causes
# the sheet's trailing filler rows to automatically extend by two rows
# (17:18) the same way the existing filler rows are shaped. Then the
# old row 6 content/format is relocated to row 8, and the old row 7
# blank-filler shape is relocated to row 9, using copy/paste-special
# (values, then formats) which reuses existing style records instead
# of minting new ones. Finally rows 5-7 are overwritten with the new
# content using styles copied from existing, equivalent rows.
# ------------------------------------------------------------------

$ws.Rows("8:9").Insert()

# Relocate old row 6 ("is_fail") content + formatting down to row 8
$ws.Range("A6:H6").Copy()
$ws.Range("A8:H8").PasteSpecial(-4122)
$ws.Range("A6:H6").Copy()
$ws.Range("A8:H8").PasteSpecial(-4163)
$ws.Rows(8).RowHeight = 35.65

# Relocate old row 7 (blank D:G filler) shape down to row 9
$ws.Range("D7:G7").Copy()
$ws.Range("D9:G9").PasteSpecial(-4122)
$ws.Range("D7:G7").Copy()
$ws.Range("D9:G9").PasteSpecial(-4163)

# ------------------------------------------------------------------
# Row 5: keep existing Arguments/SetReportErrCode row, just add the
# new docstring text into F5, and grow the row height to 24.
# ------------------------------------------------------------------
$ws.Range("F5").Value = "If iCodeBase is not iErrNotFound, calculate .iCodeReport as sum of .iCodeBase and .iCodeLocal"
$ws.Rows(5).RowHeight = 24

# ------------------------------------------------------------------
# Row 6: new "ReportError" row - copy formatting from row 5 (same
# style pattern needed: s9/s6/s3/s7/s7/s7/s7/s3) then set values.
# ------------------------------------------------------------------
$ws.Range("A5:H5").Copy()
$ws.Range("A6:H6").PasteSpecial(-4122)
$ws.Range("A6").Value = $ws.Range("A5").Value2
$ws.Range("B6").Value = "ReportError"
$ws.Range("F6").Value = "If .ErrMsg is not empty string, `nRaise error with that string as user message`nAppend .ErrMsg to .ErrMsgsAccum with added new line character if ErrMsgsAccum contains previous text"
$ws.Range("G6").Value = $ws.Range("G5").Value2
$ws.Rows(6).RowHeight = 47.25

# ------------------------------------------------------------------
# Row 7: new "ResetWarning" row - same style pattern as row 5/6.
# ------------------------------------------------------------------
$ws.Range("A5:H5").Copy()
$ws.Range("A7:H7").PasteSpecial(-4122)
$ws.Range("A7").Value = $ws.Range("A5").Value2
$ws.Range("B7").Value = "ResetWarning"
$ws.Range("F7").Value = "Re-initialize .iCodeBase, .iCodeReport, .iCodeLocal, .ErrMsg and .ErrParam class attributes to their default values as set in .__init()__"
$ws.Range("G7").Value = $ws.Range("G5").Value2
$ws.Rows(7).RowHeight = 24

$excel.CutCopyMode = 0
